$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = "'66.788.02"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -1.15%  "
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = "'3.850.91"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +2.02%  "
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.20%  "
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = "'422.33"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.77%  "
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = "'128.97"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -2.73%  "
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = "'3.847.82"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +2.54%  "
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = "'0.604"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -6.51%  "
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').Value = "'1.00"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.13%  "
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').Value = "'0.718"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -6.45%  "
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = "'0.163"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -11.00%  "
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D12').Value = "'0.0000349"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -13.36%  "
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = "'39.89"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -6.14%  "
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = "'4.462.42"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.95%  "
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('B15').Value = 'Uniswap'
$ws.Range('C15').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D15').Value = "'16.16"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +22.21%  "
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = "'9.96"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -4.18%  "
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = "'3.840.23"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +2.74%  "
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = "'0.137"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -1.62%  "
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = "'19.46"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -4.79%  "
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = "'66.954.63"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -1.05%  "
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('B21').Value = 'Polygon'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D21').Value = "'1.07"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -5.16%  "
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = "'403.24"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -8.77%  "
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').Value = "'14.20"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -9.82%  "
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = "'83.83"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -6.99%  "
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').Value = "'2.98"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -3.22%  "
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').Value = "'5.85"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +14.99%  "
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'36.81"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -3.50%  "
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = "'3.18"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -4.13%  "
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').Value = "'9.37"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -6.89%  "
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').Value = "'713.11"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +3.63%  "
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = "'2.76"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.45%  "
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = "'0.120"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -2.89%  "
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Value = "'12.27"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -2.19%  "
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').Value = "'7.36"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +3.31%  "
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = "'0.148"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -9.11%  "
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = "'37.52"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -8.85%  "
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').Value = "'0.999"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.05%  "
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').Value = "'54.69"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -5.54%  "
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = "'0.0₃0764"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +9.39%  "
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.0450"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -7.77%  "
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').Value = "'2.89"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -3.13%  "
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = "'1.00"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.58%  "
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').Value = "'0.134"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -8.77%  "
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').Value = "'4.43"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +3.32%  "
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').Value = "'3.16"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +0.15%  "
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('B46').Value = 'LidoDAOToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D46').Value = "'3.29"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -2.71%  "
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = "'143.54"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -3.30%  "
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = "'2.04"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.85%  "
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'25.75"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -6.74%  "
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = "'2.74"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -4.66%  "
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('B51').Value = 'WEMIXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D51').Value = "'2.51"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -3.31%  "
$ws.Range('E51').Style = 'Normal'
